$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(107, '2024-09-29 03:28:40', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:28:40'),
    @(108, '2024-09-29 03:28:40', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:28:40'),
    @(109, '2024-09-29 03:28:40', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:28:40'),
    @(110, '2024-09-29 03:28:41', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:28:41'),
    @(111, '2024-09-29 03:28:42', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:28:42'),
    @(112, '2024-09-29 03:32:53', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:32:53'),
    @(113, '2024-09-29 03:32:53', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:32:53'),
    @(114, '2024-09-29 03:32:53', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:32:53'),
    @(115, '2024-09-29 03:32:54', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:32:54'),
    @(116, '2024-09-29 03:32:55', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:32:55'),
    @(117, '2024-09-29 03:48:55', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:48:55'),
    @(118, '2024-09-29 03:48:55', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:48:55'),
    @(119, '2024-09-29 03:48:55', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:48:55'),
    @(120, '2024-09-29 03:48:56', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:48:56'),
    @(121, '2024-09-29 03:48:57', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:48:57'),
    @(122, '2024-09-29 03:49:42', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:49:42'),
    @(123, '2024-09-29 03:49:42', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:49:42'),
    @(124, '2024-09-29 03:49:42', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:49:42'),
    @(125, '2024-09-29 03:49:43', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:49:43'),
    @(126, '2024-09-29 03:49:44', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:49:44'),
    @(127, '2024-09-29 03:52:26', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:52:26'),
    @(128, '2024-09-29 03:52:26', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:52:26'),
    @(129, '2024-09-29 03:52:26', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:52:26'),
    @(130, '2024-09-29 03:52:27', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:52:27'),
    @(131, '2024-09-29 03:52:28', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:52:28'),
    @(132, '2024-09-29 03:53:02', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-29', '03:53:02'),
    @(133, '2024-09-29 03:53:02', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:53:02'),
    @(134, '2024-09-29 03:53:03', 'check_availability', 'https://example.com', 'Checked availability: No availability for the selected date.', '2024-09-29', '03:53:03'),
    @(135, '2024-09-29 03:53:03', 'check_availability', 'https://example.com', 'Checked availability: Selected or default date is available for booking.', '2024-09-29', '03:53:03'),
    @(136, '2024-09-29 03:53:04', 'check_availability', 'https://example.com', 'Failed to check availability: Failed to check availability', '2024-09-29', '03:53:04')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    # Column E holds a plain "YYYY-MM-DD" string; Excel would otherwise
    # auto-convert it to a date serial on input, so force Text format first
    # and then restore the default "Normal" style (matches original: no
    # explicit cell style, value stored as text).
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[5]
    $eCell.Style = "Normal"
    $ws.Cells.Item($r, 6).Value = $row[6]
}
